$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) values were updated on both the
# "展览" and "全部类型" sheets, which contain duplicated data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 22
    $ws.Range("F4").Value = 1495
    $ws.Range("F7").Value = 119
    $ws.Range("F9").Value = 285
}
